# Apply the "updates in manuscript with figures" edit to the fig1 slide.
#
# Relabels the "Categoric"/"Numeric" legend callouts to lower-case
# "discrete"/"numeric" wording (and nudges their textbox frames to match
# the new, narrower text), and swaps "categoric"/"combinations" for
# "discrete"/"combinatory" in the two annotated-matrix legends (both the
# inline callouts and the bottom key), again re-sizing the frames that
# shrank to fit the new word lengths.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$n = $s.Shapes.Count
for ($i = 1; $i -le $n; $i++) {
    $sh = $s.Shapes.Item($i)

    switch ($sh.Id) {
        604 {
            # "Categoric 1" -> "discrete 1"
            $sh.Left   = 386.93355314960627
            $sh.Top    = 385.73527559055117
            $sh.Width  = 77.8075688976378
            $sh.Height = 26.657805118110236
            $sh.TextFrame.TextRange.Text = "discrete 1"
        }
        605 {
            # "Categoric 2" -> "discrete 2"
            $sh.Left   = 415.8137696850394
            $sh.Top    = 385.73527559055117
            $sh.Width  = 77.8075688976378
            $sh.Height = 26.657805118110236
            $sh.TextFrame.TextRange.Text = "discrete 2"
        }
        606 {
            # "Numeric 1" -> "numeric 1"
            $sh.Left   = 446.26088582677164
            $sh.Top    = 387.9259251968504
            $sh.Width  = 79.92292322834646
            $sh.Height = 26.657805118110236
            $sh.TextFrame.TextRange.Text = "numeric 1"
        }
        630 {
            # "Xcategoric" -> "Xdiscrete"
            $tr = $sh.TextFrame.TextRange
            $tr.Characters(2, 9).Text = "discrete"
        }
        665 {
            # "Xcombinations" -> "Xcombinatory"
            $tr = $sh.TextFrame.TextRange
            $tr.Characters(2, 12).Text = "combinatory"
        }
        787 {
            # "Xcategoric" -> "Xdiscrete"
            $sh.Left   = 83.90875
            $sh.Top    = 577.8699114173228
            $sh.Width  = 62.64087598425197
            $sh.Height = 29.081259842519685
            $tr = $sh.TextFrame.TextRange
            $tr.Characters(2, 9).Text = "discrete"
        }
        788 {
            # "Xcombinations" -> "Xcombinatory"
            $sh.Left   = 315.84251968503935
            $sh.Top    = 577.8699114173228
            $sh.Width  = 86.34504921259843
            $sh.Height = 29.081259842519685
            $tr = $sh.TextFrame.TextRange
            $tr.Characters(2, 12).Text = "combinatory"
        }
    }
}
